# Updated cryptos list (prices + 1h volume %) per GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.575.05"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "1.584.86"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.05"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  -2.35%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -2.71%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").Value = "1.804.53"
$ws.Range("E12").Value = "  -3.04%  "
$ws.Range("D13").Value = "1.584.56"
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.529"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.42"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "26.595.54"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.62"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.72"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.40"
$ws.Range("E23").Value = "  -4.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.89"
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.42"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.43"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  -4.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.30"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.675"
$ws.Range("E33").Value = "  +24.78%  "
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").Value = "1.320.35"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0173"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.823"
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.785"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.30"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.17"
$ws.Range("E43").Value = "  -4.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.14"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").Value = "1.719.23"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.09"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.834"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0989"
$ws.Range("E49").Value = "  +3.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("E51").Value = "  -0.95%  "
